# Update gh-pages output data (苏州-漫展信息.xlsx)
# - Refresh a stale cover image URL (old image removed from CDN, replaced
#   with a re-uploaded jpeg) on both the "展览" and "全部类型" sheets.
# - Bump several "想去人数" (want-to-go count) figures to the latest scrape.

$wb = $excel.ActiveWorkbook

$oldCover = "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"
$newCover = "//i0.hdslb.com/bfs/openplatform/202406/0FW5bOPl1718591979985.jpeg"

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("I3").Value = $newCover

$ws1.Range("F8").Value  = 14754
$ws1.Range("F10").Value = 133
$ws1.Range("F11").Value = 5885
$ws1.Range("F12").Value = 602
$ws1.Range("F15").Value = 78
$ws1.Range("F17").Value = 18
$ws1.Range("F19").Value = 193
$ws1.Range("F22").Value = 94
$ws1.Range("F23").Value = 10690
$ws1.Range("F26").Value = 108
$ws1.Range("F27").Value = 3747

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("I4").Value = $newCover

$ws4.Range("F9").Value  = 14754
$ws4.Range("F11").Value = 133
$ws4.Range("F12").Value = 5886
$ws4.Range("F13").Value = 602
$ws4.Range("F16").Value = 78
$ws4.Range("F18").Value = 18
$ws4.Range("F20").Value = 193
$ws4.Range("F23").Value = 94
$ws4.Range("F25").Value = 10690
$ws4.Range("F28").Value = 108
$ws4.Range("F29").Value = 3747
